$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 435.18182
$ws.Range("I38").Value = 69.71429000000001
$ws.Range("J38").Value = 1074.75
$ws.Range("K38").Value = 209.14287
$ws.Range("L38").Value = 3224.25
$ws.Range("M38").Value = 162.85713
$ws.Range("N38").Value = -3968.25
$ws.Range("H40").Value = 1792.8572
$ws.Range("I40").Value = 1633.3334
$ws.Range("K40").Value = 1633.3334
$ws.Range("M40").Value = -1458.3334
$ws.Range("H58").Value = 3904.375
$ws.Range("I58").Value = 281.42856
$ws.Range("J58").Value = 6722.222
$ws.Range("K58").Value = 844.28568
$ws.Range("L58").Value = 20166.666
$ws.Range("M58").Value = -694.28568
$ws.Range("N58").Value = -20466.666
$ws.Range("H64").Value = 3200
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3200
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 3200
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3200
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4916
$ws.Range("H74").Value = 5224.4346
$ws.Range("I74").Value = 15333.667
$ws.Range("J74").Value = 3708.05
$ws.Range("K74").Value = 15333.667
$ws.Range("L74").Value = 3708.05
$ws.Range("M74").Value = -14397.667
$ws.Range("N74").Value = -5580.05
$ws.Range("H77").Value = 5224.4346
$ws.Range("I77").Value = 15333.667
$ws.Range("J77").Value = 3708.05
$ws.Range("K77").Value = 76668.33499999999
$ws.Range("L77").Value = 18540.25
$ws.Range("M77").Value = -71988.33499999999
$ws.Range("N77").Value = -27900.25
$ws.Range("H86").Value = 1000003
$ws.Range("I86").Value = 1000003
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1000003
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -998880
$ws.Range("N86").ClearContents()
$ws.Range("H87").Value = 27666.666
$ws.Range("J87").Value = 27375
$ws.Range("L87").Value = 27375
$ws.Range("N87").Value = -29871
$ws.Range("H89").Value = 1000003
$ws.Range("I89").Value = 1000003
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5000015
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4994399
$ws.Range("N89").ClearContents()
$ws.Range("H90").Value = 27666.666
$ws.Range("J90").Value = 27375
$ws.Range("L90").Value = 82125
$ws.Range("N90").Value = -94605

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2460.923
$ws.Range("I122").Value = 1428
$ws.Range("J122").Value = 2920
$ws.Range("K122").Value = 4284
$ws.Range("L122").Value = 8760
$ws.Range("M122").Value = -1834
$ws.Range("N122").Value = -13660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 29800
$ws.Range("J35").Value = 29800
$ws.Range("L35").Value = 29800
$ws.Range("N35").Value = -30420
$ws.Range("H82").Value = 28302.857
$ws.Range("I82").Value = 3128.5
$ws.Range("J82").Value = 38372.6
$ws.Range("K82").Value = 3128.5
$ws.Range("L82").Value = 38372.6
$ws.Range("M82").Value = -2745.5
$ws.Range("N82").Value = -39138.6
$ws.Range("H85").Value = 28302.857
$ws.Range("I85").Value = 3128.5
$ws.Range("J85").Value = 38372.6
$ws.Range("K85").Value = 3128.5
$ws.Range("L85").Value = 38372.6
$ws.Range("M85").Value = -1802.5
$ws.Range("N85").Value = -41024.6
$ws.Range("H86").Value = 2213.7144
$ws.Range("I86").Value = 2200
$ws.Range("J86").Value = 2224
$ws.Range("K86").Value = 2200
$ws.Range("L86").Value = 2224
$ws.Range("M86").Value = -1077
$ws.Range("N86").Value = -4470
$ws.Range("H89").Value = 2213.7144
$ws.Range("I89").Value = 2200
$ws.Range("J89").Value = 2224
$ws.Range("K89").Value = 11000
$ws.Range("L89").Value = 11120
$ws.Range("M89").Value = -5384
$ws.Range("N89").Value = -22352
$ws.Range("H134").Value = 23536.533
$ws.Range("I134").Value = 28049.838
$ws.Range("J134").Value = 2662.5
$ws.Range("K134").Value = 84149.514
$ws.Range("L134").Value = 7987.5
$ws.Range("M134").Value = -81614.514
$ws.Range("N134").Value = -13057.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10475
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 10475
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 10475
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -12471
$ws.Range("H83").Value = 10475
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 10475
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 52375
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -62359
$ws.Range("H102").Value = 1148.5
$ws.Range("I102").Value = 1112.7858
$ws.Range("J102").Value = 1398.5
$ws.Range("K102").Value = 1112.7858
$ws.Range("L102").Value = 1398.5
$ws.Range("M102").Value = 509.2141999999999
$ws.Range("N102").Value = -4642.5
$ws.Range("H132").Value = 185483.81
$ws.Range("I132").Value = 252790.5
$ws.Range("J132").Value = 5999.3335
$ws.Range("K132").Value = 758371.5
$ws.Range("L132").Value = 17998.0005
$ws.Range("M132").Value = -755841.5
$ws.Range("N132").Value = -23058.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 608
$ws.Range("I22").Value = 433.33334
$ws.Range("J22").Value = 695.3333
$ws.Range("K22").Value = 433.33334
$ws.Range("L22").Value = 695.3333
$ws.Range("M22").Value = -138.33334
$ws.Range("N22").Value = -1285.3333
$ws.Range("H27").Value = 608
$ws.Range("I27").Value = 433.33334
$ws.Range("J27").Value = 695.3333
$ws.Range("K27").Value = 433.33334
$ws.Range("L27").Value = 695.3333
$ws.Range("M27").Value = -326.33334
$ws.Range("N27").Value = -909.3333
$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2590
